# TC04_INS_CancerType-BroadCancerTypes.xlsx
# Cancer type facet INS - 17 test scripts
#
# The "ProgramsTab" row's TabQuery (cell B2) is updated to a new SQL query:
#  - replaces the plain "prg.data_link" column with a CASE expression that
#    falls back to prg.website when prg.data_link is populated
#  - sorts the result using LOWER(prg.program_name) instead of the raw column
# The cell keeps its wrapped text but is bumped from 11pt to 12pt font.
# The sheet's scroll position / selection is also moved up from row 5 to row 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newQuery = @'
SELECT DISTINCT 
    prg.program_name AS "Program",
    prg.website AS "Website",
    prg.focus_area AS "Focus Area",
    prg.cancer_type AS "Cancer Type",
CASE 
        WHEN prg.data_link IS NOT NULL THEN prg.website       
        ELSE prg.data_link
    END AS "Data Location Details"
FROM 
    df_program prg
WHERE 
     prg.cancer_type IN ('Broad Cancer Types')
ORDER BY 
    LOWER (prg.program_name) ASC
LIMIT 100;
'@

$cell = $ws.Range("B2")
$cell.Value = $newQuery
$cell.WrapText = $true
$cell.Font.Size = 12

# Move the visible scroll/selection up to row 3 (was row 5)
[void]$ws.Range("C3").Select()
